$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (G1, H1) with values, then copy F1's formatting onto them
# so they share the exact same style as the other header cells.
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null

# New data cells for row 2
$ws.Range("G2").Value = 0.1228190763666741
$ws.Range("H2").Value = 0.991
